# Apply updates to the "南宁-漫展信息" workbook.
# Changes affect the "展览" sheet and the "全部类型" sheet:
#   F2: 279 -> 282
#   C3: "南宁·0316全职only-全明星周末" -> "南宁·0316全职only-全明星周末（取消）"
#   G3: 65 -> "不可售" (text instead of number)
#   F4: 1067 -> 1078
#   F5/F6 (the "第一届ANE·DACG动漫嘉年华" row): 563 -> 567

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F2").Value = 282
$wsExhibition.Range("C3").Value = "南宁·0316全职only-全明星周末（取消）"
$wsExhibition.Range("G3").Value = "不可售"
$wsExhibition.Range("F4").Value = 1078
$wsExhibition.Range("F5").Value = 567

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 282
$wsAll.Range("C3").Value = "南宁·0316全职only-全明星周末（取消）"
$wsAll.Range("G3").Value = "不可售"
$wsAll.Range("F4").Value = 1078
$wsAll.Range("F6").Value = 567
